# Auto-generated edit script: updates cryptos price/volume columns (D,E) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.703.63"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "2.673.32"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'599.71"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "'167.29"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "2.672.81"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("D11").Value = "'0.158"
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  -1.69%  "
$ws.Range("D14").Value = "'27.92"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "3.164.22"
$ws.Range("E15").Value = "  -1.37%  "
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "67.661.19"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "2.676.52"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "'11.75"
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'364.23"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").Value = "'4.83"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'70.87"
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").Value = "'10.20"
$ws.Range("E27").Value = "  +2.90%  "
$ws.Range("E29").Value = "  -2.64%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'556.36"
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "'1.93"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "'0.130"
$ws.Range("E35").Value = "  -1.24%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -4.40%  "
$ws.Range("D38").Value = "'19.54"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "'155.57"
$ws.Range("E39").Value = "  -3.24%  "
$ws.Range("D40").Value = "'0.372"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'5.31"
$ws.Range("E41").Value = "  -2.16%  "
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").Value = "'17.94"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'2.52"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'40.34"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -4.85%  "
$ws.Range("D48").Value = "'0.591"
$ws.Range("E48").Value = "  -2.66%  "
$ws.Range("D49").Value = "'153.67"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -3.25%  "
